# Apply the changes described by the diff to the workbook.
$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# --- Metadata sheet updates ---
# URL
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/union-worker-indicator"
# Version
$wsMeta.Range("B3").Value = "8.0.0"
# Date
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
# Publisher
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
# Keep the "Fixed Value" for Extension.url (row 5, column Q) in sync with the new URL
$wsElem.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/union-worker-indicator"

# Clear the Constraint(s) value for the root Extension row (row 2, column AI)
$wsElem.Range("AI2").Value = ""
